$d = $word.ActiveDocument

# --- Locate the signature block's "Date:" line (the second "Date:" in the
# document -- the first one is part of the "Start Date:" heading earlier).
$dateAnchor = $d.Range(0, $d.Content.End)
$dateAnchor.Find.Execute("Date:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateAnchor = $d.Range($dateAnchor.End, $d.Content.End)
$dateAnchor.Find.Execute("Date:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$searchFrom = $dateAnchor.End

# --- First {{start_date}} on that line -> {{employer_sign_d}}, split into
# three runs: "{{" / "employer_sign_d" / "}}"
$occ1 = $d.Range($searchFrom, $d.Content.End)
$occ1.Find.Execute("{{start_date}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base1 = $occ1.Start
$occ1.Text = "{{employer_sign_d}}"

foreach ($piece in @(@(0, 2), @(2, 17))) {
    $r = $d.Range($base1 + $piece[0], $base1 + $piece[1])
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

$searchFrom = $base1 + 19

# --- Second {{start_date}} on that line -> {{ employee_sign_d}}, split into
# six runs: "{{" / " " / "employe" / "e" / "_sign_d" / "}}"
$occ2 = $d.Range($searchFrom, $d.Content.End)
$occ2.Find.Execute("{{start_date}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base2 = $occ2.Start
$occ2.Text = "{{ employee_sign_d}}"

foreach ($piece in @(@(0, 2), @(2, 3), @(3, 10), @(10, 11), @(11, 18))) {
    $r = $d.Range($base2 + $piece[0], $base2 + $piece[1])
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

Write-Host "Done"
